$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new_row[key] = old_row[value]  (1-based worksheet rows, data rows 2..31)
$mapping = @{
    2  = 23
    3  = 24
    4  = 27
    5  = 15
    6  = 16
    7  = 2
    8  = 3
    9  = 4
    10 = 25
    11 = 26
    12 = 12
    13 = 6
    14 = 7
    15 = 8
    16 = 9
    17 = 5
    18 = 29
    19 = 19
    20 = 20
    21 = 17
    22 = 18
    23 = 10
    24 = 11
    25 = 28
    26 = 30
    27 = 31
    28 = 13
    29 = 14
    30 = 21
    31 = 22
}

$firstCol = 1   # A
$lastCol  = 18  # R

# Snapshot all current values (rows 2..31, columns A..R) before overwriting anything,
# since several destination rows also act as sources.
$snapshot = @{}
for ($r = 2; $r -le 31; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c - $firstCol]
    }
}
